$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D: CURRENT_CAPACITY with per-plant values
$ws.Range("D1").Value = "CURRENT_CAPACITY"
$ws.Range("D2").Value = 80
$ws.Range("D3").Value = 50
$ws.Range("D4").Value = 50
$ws.Range("D5").Value = 100

# Widen column C slightly and size the new column D
$ws.Columns.Item(3).ColumnWidth = 29
$ws.Columns.Item(4).ColumnWidth = 17.33

# Move the selection to reflect where the user ended up after editing
$ws.Range("D8").Select()
